$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay text (preserve exact formatting)
$textCells = @("D5", "D6", "D11", "D14", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D29", "D31", "D33", "D34", "D36", "D37", "D39", "D40", "D42", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated price (D) and volume-change (E) values
$ws.Range("D2").Value = "70.957.68"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "3.848.81"
$ws.Range("E3").Value = "  +1.19%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "704.38"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").Value = "172.95"
$ws.Range("E6").Value = "  -0.54%  "
$ws.Range("D7").Value = "3.847.57"
$ws.Range("E7").Value = "  +1.23%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").Value = "7.35"
$ws.Range("E11").Value = "  -1.68%  "
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("D14").Value = "36.63"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "4.497.42"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "3.944.70"
$ws.Range("E16").Value = "  +3.71%  "
$ws.Range("D17").Value = "70.978.04"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").Value = "17.37"
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("D21").Value = "10.68"
$ws.Range("E21").Value = "  -3.98%  "
$ws.Range("D22").Value = "492.68"
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("D23").Value = "0.716"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "85.16"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").Value = "12.16"
$ws.Range("E26").Value = "  -2.30%  "
$ws.Range("D27").Value = "10.56"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("E28").Value = "  -2.41%  "
$ws.Range("D29").Value = "3.17"
$ws.Range("E29").Value = "  +4.50%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "7.50"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("E32").Value = "  -0.92%  "
$ws.Range("D33").Value = "29.44"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").Value = "0.181"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").Value = "3.804.36"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("D36").Value = "9.15"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("D39").Value = "2.37"
$ws.Range("E39").Value = "  +6.19%  "
$ws.Range("D40").Value = "6.07"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("E41").Value = "  +6.25%  "
$ws.Range("D42").Value = "3.32"
$ws.Range("E42").Value = "  -6.31%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").Value = "163.21"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").Value = "0.000306"
$ws.Range("E46").Value = "  -6.59%  "
$ws.Range("D47").Value = "48.77"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").Value = "414.81"
$ws.Range("E48").Value = "  +2.78%  "
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("D50").Value = "8.62"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").Value = "43.15"
$ws.Range("E51").Value = "  -4.25%  "
